$wb = $excel.ActiveWorkbook

# The most recent attendance sheet ("09_01") is duplicated to start a new day's
# sheet ("10_01"), inserted at the front of the workbook - just like using
# Excel's "Move or Copy... > Create a copy" on the first tab.
$src = $wb.Worksheets.Item("09_01")
$src.Copy($src)

# The freshly inserted copy lands in tab position 1 (before $src).
$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "10_01"

# Update the new sheet's header date and the two check-in/out times for the new day.
# Leading apostrophes force these to stay plain text (matching the original
# inlineStr cells) instead of being auto-parsed into date/time serials.
$newSheet.Range("B1").Value = "'10/01/20"
$newSheet.Range("C2").Value = "'18:34"
$newSheet.Range("C3").Value = "'18:36"

# A newly created sheet starts out selected at A1 - make that explicit.
$newSheet.Range("A1").Select()

# Restore the previously-active sheet ("09_01") as the active/selected tab,
# keeping its original selection untouched.
$fresh = $wb.Worksheets.Item("09_01")
$fresh.Activate()
